# Auto-generated edit script applying numeric updates to Kujata_Profits workbook
# Applies cell-value changes, additions, and removals per the target diff,
# sheet by sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 428.21054
$ws.Range("I80").Value = 238.88889
$ws.Range("J80").Value = 598.6
$ws.Range("K80").Value = 716.6666700000001
$ws.Range("L80").Value = 1795.8
$ws.Range("M80").Value = 281.3333299999999
$ws.Range("N80").Value = -3791.8
$ws.Range("H83").Value = 428.21054
$ws.Range("I83").Value = 238.88889
$ws.Range("J83").Value = 598.6
$ws.Range("K83").Value = 2150.00001
$ws.Range("L83").Value = 5387.400000000001
$ws.Range("M83").Value = 2841.99999
$ws.Range("N83").Value = -15371.4
$ws.Range("H111").Value = 4563.8
$ws.Range("J111").Value = 3500
$ws.Range("L111").Value = 10500
$ws.Range("N111").Value = -16634
$ws.Range("H113").Value = 3376.7778
$ws.Range("I113").Value = 2133.6667
$ws.Range("K113").Value = 2133.6667
$ws.Range("M113").Value = 1120.3333
$ws.Range("H116").Value = 3458.8635
$ws.Range("I116").Value = 3059.1667
$ws.Range("K116").Value = 3059.1667
$ws.Range("M116").Value = 382.8332999999998
$ws.Range("H129").Value = 878.75
$ws.Range("I129").Value = 477.875
$ws.Range("J129").Value = 915.1932
$ws.Range("K129").Value = 1433.625
$ws.Range("L129").Value = 2745.5796
$ws.Range("M129").Value = 3566.375
$ws.Range("N129").Value = -12745.5796
$ws.Range("H132").Value = 7941124
$ws.Range("I132").Value = 8337900
$ws.Range("J132").Value = 5600
$ws.Range("K132").Value = 25013700
$ws.Range("L132").Value = 16800
$ws.Range("M132").Value = -25011170
$ws.Range("N132").Value = -21860
$ws.Range("H135").Value = 750.0625
$ws.Range("I135").Value = 285.7857
$ws.Range("K135").Value = 2572.0713
$ws.Range("M135").Value = -37.07130000000006
$ws.Range("H137").Value = 1315.8
$ws.Range("I137").Value = 1216.591
$ws.Range("J137").Value = 1588.625
$ws.Range("K137").Value = 3649.773
$ws.Range("L137").Value = 4765.875
$ws.Range("M137").Value = -1099.773
$ws.Range("N137").Value = -9865.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3111.15
$ws.Range("I32").Value = 2752.3333
$ws.Range("J32").Value = 6340.5
$ws.Range("K32").Value = 2752.3333
$ws.Range("L32").Value = 6340.5
$ws.Range("M32").Value = -2465.3333
$ws.Range("N32").Value = -6914.5
$ws.Range("H33").Value = 15000
$ws.Range("I33").Value = 15000
$ws.Range("K33").Value = 15000
$ws.Range("M33").Value = -14671
$ws.Range("H74").Value = 1019.05554
$ws.Range("I74").Value = 494.08334
$ws.Range("K74").Value = 494.08334
$ws.Range("M74").Value = 379.91666
$ws.Range("H77").Value = 1019.05554
$ws.Range("I77").Value = 494.08334
$ws.Range("K77").Value = 2470.4167
$ws.Range("M77").Value = 1897.5833
$ws.Range("H88").Value = 2563
$ws.Range("I88").Value = 1860.8
$ws.Range("J88").Value = 2882.182
$ws.Range("K88").Value = 1860.8
$ws.Range("L88").Value = 2882.182
$ws.Range("M88").Value = -1454.8
$ws.Range("N88").Value = -3694.182
$ws.Range("H91").Value = 2563
$ws.Range("I91").Value = 1860.8
$ws.Range("J91").Value = 2882.182
$ws.Range("K91").Value = 1860.8
$ws.Range("L91").Value = 2882.182
$ws.Range("M91").Value = -456.8
$ws.Range("N91").Value = -5690.182
$ws.Range("H110").Value = 1345.7646
$ws.Range("I110").Value = 812.6923
$ws.Range("K110").Value = 812.6923
$ws.Range("M110").Value = 1232.3077

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3985.92
$ws.Range("I86").Value = 4323.2104
$ws.Range("J86").Value = 2917.8333
$ws.Range("K86").Value = 4323.2104
$ws.Range("L86").Value = 2917.8333
$ws.Range("M86").Value = -3200.2104
$ws.Range("N86").Value = -5163.8333
$ws.Range("H89").Value = 3985.92
$ws.Range("I89").Value = 4323.2104
$ws.Range("J89").Value = 2917.8333
$ws.Range("K89").Value = 21616.052
$ws.Range("L89").Value = 14589.1665
$ws.Range("M89").Value = -16000.052
$ws.Range("N89").Value = -25821.1665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1016.0526
$ws.Range("I31").Value = 961.44446
$ws.Range("J31").Value = 1999
$ws.Range("K31").Value = 961.44446
$ws.Range("L31").Value = 1999
$ws.Range("M31").Value = -666.44446
$ws.Range("N31").Value = -2589
$ws.Range("H34").Value = 1016.0526
$ws.Range("I34").Value = 961.44446
$ws.Range("J34").Value = 1999
$ws.Range("K34").Value = 961.44446
$ws.Range("L34").Value = 1999
$ws.Range("M34").Value = -759.44446
$ws.Range("N34").Value = -2403
$ws.Range("H132").Value = 5217.3335
$ws.Range("I132").Value = 6569.7
$ws.Range("K132").Value = 19709.1
$ws.Range("M132").Value = -17179.1
$ws.Range("H134").Value = 9805238
$ws.Range("I134").Value = 11112422
$ws.Range("K134").Value = 33337266
$ws.Range("M134").Value = -33334731

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 639.7273
$ws.Range("I107").Value = 263.33334
$ws.Range("J107").Value = 780.875
$ws.Range("K107").Value = 790.0000200000001
$ws.Range("L107").Value = 2342.625
$ws.Range("M107").Value = 1129.99998
$ws.Range("N107").Value = -6182.625
$ws.Range("H131").Value = 11765927
$ws.Range("I131").Value = 333333570
$ws.Range("K131").Value = 1000000710
$ws.Range("M131").Value = -999995670

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32149158
$ws.Range("I70").Value = 27783578
$ws.Range("K70").Value = 27783578
$ws.Range("M70").Value = -27783308
$ws.Range("H73").Value = 32149158
$ws.Range("I73").Value = 27783578
$ws.Range("K73").Value = 27783578
$ws.Range("M73").Value = -27782642

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H124:L124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:L134").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:L139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5634
$ws.Range("M32").ClearContents()
$ws.Range("H96").Value = 1304.875
$ws.Range("I96").Value = 746.5
$ws.Range("K96").Value = 746.5
$ws.Range("M96").Value = 626.5
$ws.Range("H100").Value = 2010.4286
$ws.Range("I100").Value = 1014.8
$ws.Range("J100").Value = 4499.5
$ws.Range("K100").Value = 2029.6
$ws.Range("L100").Value = 8999
$ws.Range("M100").Value = -1488.6
$ws.Range("N100").Value = -10081
$ws.Range("H140").Value = 32056.334
$ws.Range("J140").Value = 32056.334
$ws.Range("L140").Value = 32056.334
$ws.Range("N140").Value = -42416.334
